# Add the missing "Q3" pick-and-place row.
# A new row is inserted above the current row 2 (pushing all existing
# data rows down by one), and populated with the Q3 component's data:
#   Designator = Q3, Mid X = 38.29, Mid Y = 89.1, Layer = Top, Rotation = 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2; everything below (rows 2:37) shifts to 3:38.
$ws.Rows("2:2").Insert()

# The inserted row inherits formatting from the header row above it by
# default; strip that so the new data row matches the rest of the table
# (no explicit cell style).
$ws.Range("A2:E2").ClearFormats()

# Populate the new Q3 row.
$ws.Range("A2").Value = "Q3"
$ws.Range("B2").Value = 38.29
$ws.Range("C2").Value = 89.1
$ws.Range("D2").Value = "Top"
$ws.Range("E2").Value = 0

# Match the saved cursor position from the edit.
[void]$ws.Range("D6").Select()
